# Add 6 new LeetCode problem rows (142-147) to the tracking sheet, plus
# three blank placeholder rows (139-141) that only carry date-column
# formatting (mirrors the pattern already used elsewhere in the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy the date-format (H/I, style index 4) down for the three blank
#     "spacer" rows, matching the existing H138:I138 formatting. ---
$ws.Range("H138:I138").Copy()
$ws.Range("H139:I139").PasteSpecial(-4122)
$ws.Range("H140:I140").PasteSpecial(-4122)
$ws.Range("H141:I141").PasteSpecial(-4122)

# Row 142: 498 - Diagonal Traverse
$ws.Range("A142").Value = 498
$ws.Range("B142").Value = "Diagonal Traverse"
$ws.Range("C142").Value = "#array #matrix "
$ws.Range("D142").Value = "medium"
$ws.Range("E142").Value = 1
$ws.Range("F142").Value = 0
$ws.Range("G142").Value = 25
$ws.Range("H138:I138").Copy()
$ws.Range("H142:I142").PasteSpecial(-4122)
$ws.Range("H142").Value = 45894
$ws.Range("I142").Value = 45894
$ws.Rows.Item(142).RowHeight = 17

# Row 143: 3000 - Maximum Area of Longest Diagonal Rectangle
$ws.Range("A143").Value = 3000
$ws.Range("B143").Value = "Maximum Area of Longest Diagonal Rectangle"
$ws.Range("D143").Value = "easy"
$ws.Range("E143").Value = 1
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 5
$ws.Range("H138:I138").Copy()
$ws.Range("H143:I143").PasteSpecial(-4122)
$ws.Range("H143").Value = 45895
$ws.Range("I143").Value = 45895
$ws.Rows.Item(143).RowHeight = 51

# Row 144: 1181 - Before and After Puzzle
$ws.Range("A144").Value = 1181
$ws.Range("B144").Value = "Before and After Puzzle"
$ws.Range("C144").Value = "#array #hash-table #string #sorting"
$ws.Range("D144").Value = "medium"
$ws.Range("E144").Value = 0
$ws.Range("F144").Value = 1
$ws.Range("G144").Value = 25
$ws.Range("H138:I138").Copy()
$ws.Range("H144:I144").PasteSpecial(-4122)
$ws.Range("H144").Value = 45895
$ws.Range("I144").Value = 45895
$ws.Rows.Item(144).RowHeight = 34

# Row 145: 3459 - Length of Longest V-Shaped Diagonal Segment
$ws.Range("A145").Value = 3459
$ws.Range("B145").Value = "Length of Longest V-Shaped Diagonal Segment"
$ws.Range("C145").Value = "#matrix #dfs #cache"
$ws.Range("D145").Value = "hard"
$ws.Range("E145").Value = 0
$ws.Range("F145").Value = 1
$ws.Range("G145").Value = 30
$ws.Range("H138:I138").Copy()
$ws.Range("H145:I145").PasteSpecial(-4122)
$ws.Range("H145").Value = 45896
$ws.Range("I145").Value = 45896
$ws.Range("J145").Value = "@cache?"
$ws.Rows.Item(145).RowHeight = 51

# Row 146: 3446 - Sort Matrix by Diagonals
$ws.Range("A146").Value = 3446
$ws.Range("B146").Value = "Sort Matrix by Diagonals"
$ws.Range("C146").Value = "#matrix #sorting "
$ws.Range("D146").Value = "medium"
$ws.Range("E146").Value = 1
$ws.Range("F146").Value = 0
$ws.Range("G146").Value = 22
$ws.Range("H138:I138").Copy()
$ws.Range("H146:I146").PasteSpecial(-4122)
$ws.Range("H146").Value = 45897
$ws.Range("I146").Value = 45897
$ws.Rows.Item(146).RowHeight = 17

# Row 147: 3021 - Alice and Bob Playing Flower Game
$ws.Range("A147").Value = 3021
$ws.Range("B147").Value = "Alice and Bob Playing Flower Game"
$ws.Range("C147").Value = "#math "
$ws.Range("D147").Value = "medium"
$ws.Range("E147").Value = 1
$ws.Range("F147").Value = 0
$ws.Range("G147").Value = 10
$ws.Range("H138:I138").Copy()
$ws.Range("H147:I147").PasteSpecial(-4122)
$ws.Range("H147").Value = 45898
$ws.Range("I147").Value = 45898
$ws.Rows.Item(147).RowHeight = 34

$ws.Range("J147").Select()
